$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before current row 1279 (i.e. after the two
# "1a (guarda)" rows for 2023-09-04), shifting all subsequent rows down by 2.
$ws.Rows.Item(1279).Insert()
$ws.Rows.Item(1279).Insert()

# Fill the first new row (1279) with a new weekly price observation.
$ws.Cells.Item(1279, 1).Value = 7
$ws.Cells.Item(1279, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(1279, 3).Value = "Ñuble"
$ws.Cells.Item(1279, 4).Value = 45239
$ws.Cells.Item(1279, 5).Value = 16
$ws.Cells.Item(1279, 6).Value = 100112004
$ws.Cells.Item(1279, 7).Value = "Cebolla"
$ws.Cells.Item(1279, 8).Value = "Sin especificar"
$ws.Cells.Item(1279, 9).Value = "1a nueva(o)"
$ws.Cells.Item(1279, 10).Value = 20000
$ws.Cells.Item(1279, 11).Value = 3500
$ws.Cells.Item(1279, 12).Value = 3500
$ws.Cells.Item(1279, 13).Value = 3500
$ws.Cells.Item(1279, 14).Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(1279, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1279, 16).Value = 350
$ws.Cells.Item(1279, 17).Value = 10
$ws.Cells.Item(1279, 18).Value = "Hortaliza"

# Fill the second new row (1280) with a new weekly price observation.
$ws.Cells.Item(1280, 1).Value = 7
$ws.Cells.Item(1280, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(1280, 3).Value = "Ñuble"
$ws.Cells.Item(1280, 4).Value = 45239
$ws.Cells.Item(1280, 5).Value = 16
$ws.Cells.Item(1280, 6).Value = 100112004
$ws.Cells.Item(1280, 7).Value = "Cebolla"
$ws.Cells.Item(1280, 8).Value = "Sin especificar"
$ws.Cells.Item(1280, 9).Value = "2a nueva(o)"
$ws.Cells.Item(1280, 10).Value = 25000
$ws.Cells.Item(1280, 11).Value = 3000
$ws.Cells.Item(1280, 12).Value = 3000
$ws.Cells.Item(1280, 13).Value = 3000
$ws.Cells.Item(1280, 14).Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(1280, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1280, 16).Value = 300
$ws.Cells.Item(1280, 17).Value = 10
$ws.Cells.Item(1280, 18).Value = "Hortaliza"
